$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 108; existing rows 108:130 shift down to 110:132
$ws.Rows("108:109").Insert()

# --- New row 108 ---
$ws.Cells.Item(108, 1).Value = 5
$ws.Cells.Item(108, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(108, 3).Value = "Maule"
$ws.Cells.Item(108, 4).Value = 44932
$ws.Cells.Item(108, 5).Value = 7
$ws.Cells.Item(108, 6).Value = "Fruta"
$ws.Cells.Item(108, 7).Value = 100103
$ws.Cells.Item(108, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(108, 9).Value = 100103002
$ws.Cells.Item(108, 10).Value = "Ciruela"
$ws.Cells.Item(108, 11).Value = "Black Amber"
$ws.Cells.Item(108, 12).Value = "Primera"
$ws.Cells.Item(108, 13).Value = 180
$ws.Cells.Item(108, 14).Value = 15000
$ws.Cells.Item(108, 15).Value = 15000
$ws.Cells.Item(108, 16).Value = 15000
$ws.Cells.Item(108, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(108, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(108, 19).Value = 833
$ws.Cells.Item(108, 20).Value = 18

# --- New row 109 ---
$ws.Cells.Item(109, 1).Value = 5
$ws.Cells.Item(109, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(109, 3).Value = "Maule"
$ws.Cells.Item(109, 4).Value = 44932
$ws.Cells.Item(109, 5).Value = 7
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100103
$ws.Cells.Item(109, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(109, 9).Value = 100103002
$ws.Cells.Item(109, 10).Value = "Ciruela"
$ws.Cells.Item(109, 11).Value = "Black Amber"
$ws.Cells.Item(109, 12).Value = "Segunda"
$ws.Cells.Item(109, 13).Value = 150
$ws.Cells.Item(109, 14).Value = 12000
$ws.Cells.Item(109, 15).Value = 12000
$ws.Cells.Item(109, 16).Value = 12000
$ws.Cells.Item(109, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(109, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(109, 19).Value = 667
$ws.Cells.Item(109, 20).Value = 18
